# Apply the "Add files via upload" edit to 8-2-1.xlsx:
#  - Update the organization website URL from "www.stat.kg" to "www.stat.gov.kg"
#  - Turn on word-wrap for the Goal cell (B2) to match the wrapped style used
#    by the other long text cells in column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Organization website (if available) -> B10
$ws.Range("B10").Value = "www.stat.gov.kg"

# Goal cell (B2): switch its cell style to the word-wrap variant
$ws.Range("B2").WrapText = $true
